# storm2_results.xlsx edit script
#
# Commit intent ("processing to get transported sizes!!!!!!"): the CP:CS
# columns held timestamps/values "transported" in from a different storm's
# run (2021-07-23) that didn't belong in this 2022-08-03 aggregate table, so
# their contents are cleared out. The bottom summary block (rows 34-60,
# one row per metric in A2:A28) is then extended down to also cover the
# "tau" / "amp_factor" rows (29-30) that previously had no summary rows,
# and a few number formats in that block are cleaned up along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clear out the stray CP:CS data (columns with mismatched timestamps from a
#    different storm/date range) across the whole data block, rows 1-31.
# ---------------------------------------------------------------------------
$ws.Range("CP1:CS1").ClearContents() | Out-Null
$ws.Range("CP2:CS30").ClearContents() | Out-Null
$ws.Range("CP31:CS31").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 2. A few of the "max" helper formulas in the summary block were manually
#    tightened to stop at CO (the real last data column) instead of CS.
# ---------------------------------------------------------------------------
$ws.Range("B40").Formula = "=MAX(B8:CO8)"
$ws.Range("B41").Formula = "=MAX(B9:CO9)"
$ws.Range("B42").Formula = "=MAX(B10:CO10)"

# ---------------------------------------------------------------------------
# 3. Extend the summary block down two more rows (61, 62) to cover the "tau"
#    (row 29) and "amp_factor" (row 30) series, same A:C pattern as every
#    other row in the block, then leave row 63 blank as a spacer (matching
#    the look of row 33 above the block).
# ---------------------------------------------------------------------------
$ws.Range("A61").Value = "tau"
$ws.Range("B61").Formula = "=MAX(B29:CS29)"
$ws.Range("C61").Formula = "=_xlfn.XLOOKUP(MAX(B29:CS29), B29:CS29, `$B`$1:`$CS`$1)"

$ws.Range("A62").Value = "amp_factor"
$ws.Range("B62").Formula = "=MAX(B30:CS30)"
$ws.Range("C62").Formula = "=_xlfn.XLOOKUP(MAX(B30:CS30), B30:CS30, `$B`$1:`$CS`$1)"

# ---------------------------------------------------------------------------
# 4. Formatting touch-ups in the B column of the summary block: the
#    scientific-notation format is dropped everywhere. The amplitude rows
#    get a fixed 6-decimal format, the grey-highlighted "channel 1-3" rows
#    fall back to General (keeping their grey fill), and the rest fall back
#    to plain General with no fill.
# ---------------------------------------------------------------------------
$ws.Range("B34:B36").NumberFormat = "0.000000"

$ws.Range("B37:B39").NumberFormat = "General"
$ws.Range("B43:B51").NumberFormat = "General"

$ws.Range("B40:B42").NumberFormat = "General"
$ws.Range("B52:B60").NumberFormat = "General"

$ws.Range("B61:B62").NumberFormat = "General"
$ws.Range("B61:B62").Interior.Color = $ws.Range("B60").Interior.Color
$ws.Range("B61:B62").HorizontalAlignment = -4108

$ws.Range("C61:C62").NumberFormat = $ws.Range("C60").NumberFormat
$ws.Range("C61:C62").Interior.Color = $ws.Range("C60").Interior.Color

$ws.Range("B63").NumberFormat = "m/d/yy h:mm"
$ws.Range("B63").HorizontalAlignment = -4108
$ws.Range("C63").NumberFormat = "[$-F400]h:mm:ss AM/PM"

# ---------------------------------------------------------------------------
# 5. Column widths: B and C best-fit again now that their displayed content
#    changed (longer "0.000000" numbers / different date text).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.43
$ws.Columns.Item(3).ColumnWidth = 11.29
$ws.Range("B1:C1").EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 6. View state: scroll position / active selection on the sheet.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("Z41").Select() | Out-Null
